$d = $word.ActiveDocument

# 1. Skills line: insert "HTML, PHP, CSS, " before "MatLab, Verilog, "
$d.Content.Find.Execute("SQL, MatLab, Verilog,", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "SQL, HTML, PHP, CSS, MatLab, Verilog,", 2)

# 2. Skills line: insert "Web Development, " after "App Development, " before "Unity, Simulink, Solidworks"
$d.Content.Find.Execute("App Development, Unity, Simulink, Solidworks", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "App Development, Web Development, Unity, Simulink, Solidworks", 2)

Write-Host "Done"
